$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Values for columns E (index) and F (value) for rows 20-31
$data = @{
    20 = @(19, 12)
    21 = @(20, 11)
    22 = @(21, 11)
    23 = @(22, 12)
    24 = @(23, 12)
    25 = @(24, 11)
    26 = @(25, 14)
    27 = @(26, 14)
    28 = @(27, 13)
    29 = @(28, 12)
    30 = @(29, 11)
    31 = @(30, 11)
}

foreach ($row in 20..31) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 5).Value = $vals[0]
    $ws.Cells.Item($row, 6).Value = $vals[1]
}

# Copy the formatting from an existing styled E/F cell (row 19) onto the new cells
$ws.Range("E19:F19").Copy()
$ws.Range("E20:F31").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
